{"js": "// Update the date heading and the 20x5 grid of arithmetic problems.\n// The replacement grid values, row-major (top-to-bottom, left-to-right),\n// matching the table's existing row/column layout exactly.\nconst newTableValues = [\n  [\"72-68=\", \"6+45=\", \"64-31=\", \"31+64=\", \"77-75=\"],\n  [\"46+3=\", \"27+66=\", \"40+10=\", \"11+57=\", \"50+44=\"],\n  [\"15+2=\", \"60+36=\", \"42+9=\", \"34+39=\", \"99-14=\"],\n  [\"74-17=\", \"70-60=\", \"56-55=\", \"84-43=\", \"0+14=\"],\n  [\"35+59=\", \"71+9=\", \"63-23=\", \"79-4=\", \"46-38=\"],\n  [\"75-6=\", \"86-73=\", \"13+4=\", \"32+47=\", \"15+26=\"],\n  [\"47+13=\", \"83+0=\", \"19+39=\", \"94-78=\", \"79-54=\"],\n  [\"40-11=\", \"59-6=\", \"15-0=\", \"53+13=\", \"24+68=\"],\n  [\"23+7=\", \"20+47=\", \"86-53=\", \"37+41=\", \"67-65=\"],\n  [\"62-16=\", \"73-57=\", \"93-49=\", \"99-60=\", \"42-32=\"],\n  [\"52-32=\", \"67+23=\", \"71-23=\", \"10+22=\", \"40+44=\"],\n  [\"73-72=\", \"73-31=\", \"52-47=\", \"5+52=\", \"39+2=\"],\n  [\"87+0=\", \"0+17=\", \"46-32=\", \"6+29=\", \"90-56=\"],\n  [\"82-25=\", \"81-40=\", \"91+8=\", \"71+10=\", \"57-49=\"],\n  [\"99-53=\", \"65-27=\", \"1+43=\", \"11+37=\", \"37+12=\"],\n  [\"38-35=\", \"85-67=\", \"21+36=\", \"78-46=\", \"26+6=\"],\n  [\"48+15=\", \"36+25=\", \"43-37=\", \"52+26=\", \"2+66=\"],\n  [\"97-14=\", \"22+3=\", \"21+42=\", \"75-51=\", \"43+23=\"],\n  [\"74+10=\", \"31-2=\", \"14+33=\", \"12+49=\", \"18+1=\"],\n  [\"8+45=\", \"6+84=\", \"3+60=\", \"72+22=\", \"68+26=\"],\n];\n\n// 1. Update the title paragraph (\"2023-10-15 Sunday\" -> \"2023-10-16 Monday\"),\n// replacing only the text run content so the existing run formatting\n// (Arial, size 30 half-points) stays untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text === \"2023-10-15 Sunday\") {\n    p.insertText(\"2023-10-16 Monday\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2. Update the table of math problems in place, preserving per-cell\n// formatting (font/size) by only changing the text value of each cell.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newTableValues;\nawait context.sync();\n", "ps1": "# Update the date heading and the 20x5 grid of arithmetic problems.\n$d = $word.ActiveDocument\n\n# 1. Update the title paragraph (\"2023-10-15 Sunday\" -> \"2023-10-16 Monday\").\n# Setting Range.Text in place keeps the existing run formatting\n# (Arial, size 30 half-points) untouched.\n$d.Paragraphs(1).Range.Text = \"2023-10-16 Monday\"\n\n# 2. Update the table of math problems, cell by cell (row-major order,\n# matching the existing table layout exactly) so the duplicate value\n# \"71-66=\" (rows 2 and 18, column 2) is replaced positionally rather than\n# via a global text replace.\n$newValues = @(\n    @(\"72-68=\", \"6+45=\", \"64-31=\", \"31+64=\", \"77-75=\"),\n    @(\"46+3=\", \"27+66=\", \"40+10=\", \"11+57=\", \"50+44=\"),\n    @(\"15+2=\", \"60+36=\", \"42+9=\", \"34+39=\", \"99-14=\"),\n    @(\"74-17=\", \"70-60=\", \"56-55=\", \"84-43=\", \"0+14=\"),\n    @(\"35+59=\", \"71+9=\", \"63-23=\", \"79-4=\", \"46-38=\"),\n    @(\"75-6=\", \"86-73=\", \"13+4=\", \"32+47=\", \"15+26=\"),\n    @(\"47+13=\", \"83+0=\", \"19+39=\", \"94-78=\", \"79-54=\"),\n    @(\"40-11=\", \"59-6=\", \"15-0=\", \"53+13=\", \"24+68=\"),\n    @(\"23+7=\", \"20+47=\", \"86-53=\", \"37+41=\", \"67-65=\"),\n    @(\"62-16=\", \"73-57=\", \"93-49=\", \"99-60=\", \"42-32=\"),\n    @(\"52-32=\", \"67+23=\", \"71-23=\", \"10+22=\", \"40+44=\"),\n    @(\"73-72=\", \"73-31=\", \"52-47=\", \"5+52=\", \"39+2=\"),\n    @(\"87+0=\", \"0+17=\", \"46-32=\", \"6+29=\", \"90-56=\"),\n    @(\"82-25=\", \"81-40=\", \"91+8=\", \"71+10=\", \"57-49=\"),\n    @(\"99-53=\", \"65-27=\", \"1+43=\", \"11+37=\", \"37+12=\"),\n    @(\"38-35=\", \"85-67=\", \"21+36=\", \"78-46=\", \"26+6=\"),\n    @(\"48+15=\", \"36+25=\", \"43-37=\", \"52+26=\", \"2+66=\"),\n    @(\"97-14=\", \"22+3=\", \"21+42=\", \"75-51=\", \"43+23=\"),\n    @(\"74+10=\", \"31-2=\", \"14+33=\", \"12+49=\", \"18+1=\"),\n    @(\"8+45=\", \"6+84=\", \"3+60=\", \"72+22=\", \"68+26=\")\n)\n\n$t = $d.Tables(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
